# Fix Training Data Issue (#48)
#
# The "Date" column (BF) on Sheet1 was stamped with the literal text
# "6-15-2011-12" for every stat row. Because of the way the NBA stats
# export showed dates, the data was off by a day relative to the actual
# game date. Correct it to the proper ISO date string "2012-06-15" for
# every data row (BF2:BF31), keeping the column as plain text so Excel
# does not reinterpret it as a date serial value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 31
$dateCol  = "BF"
$correctedDate = "2012-06-15"

# Make sure the column keeps storing a literal string instead of being
# auto-converted to a date serial by Excel's input parsing.
$ws.Range("$dateCol$firstRow`:$dateCol$lastRow").NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Range("$dateCol$row").Value = $correctedDate
}
